# Update the "想去人数" (number of people wanting to go) figures in the
# F column across all four sheets, per the latest gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2218
$ws.Range("F5").Value = 4366
$ws.Range("F13").Value = 680251
$ws.Range("F14").Value = 1683
$ws.Range("F19").Value = 1310
$ws.Range("F21").Value = 1167
$ws.Range("F23").Value = 1585
$ws.Range("F33").Value = 592
$ws.Range("F35").Value = 3173
$ws.Range("F40").Value = 2653
$ws.Range("F44").Value = 1031

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 145600
$ws.Range("F10").Value = 145600
$ws.Range("F19").Value = 429
$ws.Range("F24").Value = 676

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3143
$ws.Range("F9").Value = 653
$ws.Range("F11").Value = 158
$ws.Range("F12").Value = 2065

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 653
$ws.Range("F7").Value = 2218
$ws.Range("F8").Value = 158
$ws.Range("F9").Value = 2065
$ws.Range("F10").Value = 4366
$ws.Range("F16").Value = 680251
$ws.Range("F19").Value = 1683
$ws.Range("F20").Value = 145600
$ws.Range("F24").Value = 1310
$ws.Range("F26").Value = 1167
$ws.Range("F28").Value = 1585
$ws.Range("F40").Value = 3173
$ws.Range("F46").Value = 2653
$ws.Range("F51").Value = 1031
